$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" row correct-answer marks (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update the "Total" row correct-answer marks (B12): 15 -> 25
$ws.Range("B12").Value = 25

# Update the correct/total marks label (E12): "3/84" -> "25/140"
$ws.Range("E12").Value = "25/140"
